$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.201663017272949
$ws.Range("B1").Value = 5.601093292236328
$ws.Range("C1").Value = 4.481451511383057
$ws.Range("D1").Value = 5.155181884765625
$ws.Range("E1").Value = 5.543619632720947
